# "connected to mongodb cloud"
# Adds a new "Phase 5" slide (slide33) at the end of the deck, built by
# duplicating the previous "Phase 5" slide (slide32) and then:
#   - removing the "404.Ejs" screenshot + its caption textbox
#   - repositioning/rewriting the remaining caption textbox with the
#     MongoDB-cloud-connection blurb
# All other shapes (title box, db.js code screenshot, highlight rectangles,
# terminal screenshot, bent arrow) are carried over unchanged from slide32.

$p = $ppt.ActivePresentation

# --- duplicate the last "Phase 5" slide (slide 32) to create slide 33 ---
$src = $p.Slides.Item($p.Slides.Count)
$dup = $src.Duplicate()
$s = $dup.Item(1)

# --- drop the 404.ejs screenshot and its caption ---
$s.Shapes.Item("Picture 2").Delete()
$s.Shapes.Item("TextBox 14").Delete()

# --- rewrite + reposition the remaining caption textbox ---
$tb = $s.Shapes.Item("TextBox 16")
$tr = $tb.TextFrame.TextRange
$tr.Text = "Connecting to MongoDB online is so simple and straight forward:" + "`r" + "`r" + "Just sign up, create user, get link and use in the app" + "`r" + "`r" + "`r" + "Heroku deployment is as usual."

$tb.Left = 12.268425496850394
$tb.Top = 63.789213198425195
$tb.Width = 487.4896087992126
$tb.Height = 123.59527509055118
